# Adds a new "2020" column (column L) to the SDG 3.5.2.1 indicator table,
# mirroring the formatting of the existing year columns, and moves the
# active selection to J6 (as recorded in the author's last saved view).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 2 (blank cell under the thick rule, right of K2). K2 carries the
# thick-bottom border used across the row; the new cell additionally
# switches to the plain (non-bold, 11pt) Times New Roman font used by
# the sheet's base style.
# ---------------------------------------------------------------------
$ws.Range("K2").Copy()
$ws.Range("L2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("L2").Font.Name = "Times New Roman"
$ws.Range("L2").Font.Size = 11
$ws.Range("L2").Font.Bold = $false

# ---------------------------------------------------------------------
# Row 3: year header "2020", formatted like the other year cells (K3).
# ---------------------------------------------------------------------
$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)
$ws.Range("L3").Value = 2020

# ---------------------------------------------------------------------
# Row 4: first data row of the table (bold indicator total), formatted
# like K4.
# ---------------------------------------------------------------------
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L4").Value = 15.121199070053857

# ---------------------------------------------------------------------
# Rows 5-32: regular data rows, each formatted like its row's existing
# K-column cell (style is identical - fontId 3, no border - for all of
# them), then filled with the 2020 values.
# ---------------------------------------------------------------------
$values = @(
  4.0418020313117182,
  26.038712293651503,
  6.6322037707763002,
  0.72337439688659655,
  12.766404830206815,
  3.4408726052927023,
  0.31797968427797146,
  6.6053177641160472,
  3.8076839061225556,
  0.40244525738386439,
  7.1853419025188616,
  21.314558187024932,
  6.076523688314845,
  37.122644813336137,
  17.175755080979702,
  1.8631795589137379,
  32.838796692664047,
  26.748199140342599,
  8.819184666377593,
  45.075501464953796,
  22.63720315804489,
  6.5075739764215905,
  38.284059576513137,
  23.966278037020849,
  11.409880956908683,
  35.079762166299332,
  16.385848928775125,
  3.2448358437546645
)

for ($i = 0; $i -lt $values.Count; $i++) {
  $row = 5 + $i
  $ws.Range("K$row").Copy()
  $ws.Range("L$row").PasteSpecial(-4122)
  $ws.Range("L$row").Value = $values[$i]
}

# ---------------------------------------------------------------------
# Row 33: last row of the table, above the thick closing rule -
# formatted like K33.
# ---------------------------------------------------------------------
$ws.Range("K33").Copy()
$ws.Range("L33").PasteSpecial(-4122)
$ws.Range("L33").Value = 28.789140981035917

# ---------------------------------------------------------------------
# Restore the selection recorded in the workbook's last saved view.
# ---------------------------------------------------------------------
$ws.Range("J6").Select()
